# Commit: "updating the new files"
# - Adds two new worksheets (Gmail_Test_1, Gmail_Test_2) with the same
#   login-style two-row/two-column layout used by the existing sheets,
#   each cell text hyperlinked to itself via a mailto: link (matching the
#   pattern already used on login_Test_Case_01).
# - Moves the active tab from Google_Search_Test_2 (index 2) to the newly
#   added Gmail_Test_2 (index 4), which becomes the selected tab.
# - Updates the selection left behind on Google_Search_Test_2 since it is
#   no longer the active sheet.

$wb = $excel.ActiveWorkbook

# --- Google_Search_Test_2 (3rd tab) loses "active" status; its lingering
#     selection moves to B8 and tabSelected is dropped (handled implicitly
#     once a different sheet becomes active below). ---
$googleTest2 = $wb.Worksheets.Item(3)
$googleTest2.Range("B8").Select()

# --- Gmail_Test_1 -----------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$gmail1 = $wb.Worksheets.Add($null, $lastSheet)
$gmail1.Name = "Gmail_Test_1"

$gmail1.Range("A1").Value = "Username"
$gmail1.Range("B1").Value = "password"
$gmail1.Range("A2").Value = "venkatatest1729@gmail.com"
$gmail1.Range("B2").Value = "Infy@123"

$gmail1.Hyperlinks.Add($gmail1.Range("A2"), "mailto:venkatatest1729@gmail.com")
$gmail1.Hyperlinks.Add($gmail1.Range("B2"), "mailto:Infy@123")

$gmail1.Columns.Item(1).ColumnWidth = 22.736979166666668
$gmail1.Columns.Item(2).ColumnWidth = 9.166666666666666

$gmail1.Range("A1:XFD1048576").Select()

# --- Gmail_Test_2 -------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$gmail2 = $wb.Worksheets.Add($null, $lastSheet2)
$gmail2.Name = "Gmail_Test_2"

$gmail2.Range("A1").Value = "Username"
$gmail2.Range("B1").Value = "password"
$gmail2.Range("A2").Value = "venkatatest1729@gmail.com"
$gmail2.Range("B2").Value = "Infy@123"

$gmail2.Hyperlinks.Add($gmail2.Range("A2"), "mailto:venkatatest1729@gmail.com")
$gmail2.Hyperlinks.Add($gmail2.Range("B2"), "mailto:Infy@123")

$gmail2.Columns.Item(1).ColumnWidth = 22.736979166666668
$gmail2.Columns.Item(2).ColumnWidth = 9.166666666666666

$gmail2.Range("J10").Select()

# Gmail_Test_2 is the new active / selected tab.
$gmail2.Activate()
